$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 280.33334
$ws.Range("I6").Value = 280.33334
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 841.0000200000001
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -729.0000200000001
$ws.Range("N6").Value = ""

$ws.Range("H19").Value = 42721.2
$ws.Range("I19").Value = 633
$ws.Range("J19").Value = 60759
$ws.Range("K19").Value = 633
$ws.Range("L19").Value = 60759
$ws.Range("M19").Value = -458
$ws.Range("N19").Value = -61109

$ws.Range("H43").Value = 7639.8
$ws.Range("J43").Value = 8987.5
$ws.Range("L43").Value = 8987.5
$ws.Range("N43").Value = -9125.5

$ws.Range("H70").Value = 841.5833
$ws.Range("I70").Value = 758.3333
$ws.Range("J70").Value = 924.8333
$ws.Range("K70").Value = 2274.9999
$ws.Range("L70").Value = 2774.4999
$ws.Range("M70").Value = -2004.9999
$ws.Range("N70").Value = -3314.4999

$ws.Range("H73").Value = 841.5833
$ws.Range("I73").Value = 758.3333
$ws.Range("J73").Value = 924.8333
$ws.Range("K73").Value = 2274.9999
$ws.Range("L73").Value = 2774.4999
$ws.Range("M73").Value = -1338.9999
$ws.Range("N73").Value = -4646.4999

$ws.Range("H80").Value = 1895.2
$ws.Range("I80").Value = 2938.5
$ws.Range("J80").Value = 1199.6666
$ws.Range("K80").Value = 8815.5
$ws.Range("L80").Value = 3598.9998
$ws.Range("M80").Value = -7817.5
$ws.Range("N80").Value = -5594.9998

$ws.Range("H83").Value = 1895.2
$ws.Range("I83").Value = 2938.5
$ws.Range("J83").Value = 1199.6666
$ws.Range("K83").Value = 26446.5
$ws.Range("L83").Value = 10796.9994
$ws.Range("M83").Value = -21454.5
$ws.Range("N83").Value = -20780.9994

$ws.Range("H98").Value = 1673.8334
$ws.Range("I98").Value = 1780.4546
$ws.Range("K98").Value = 1780.4546
$ws.Range("M98").Value = -282.4546

$ws.Range("H107").Value = 1059.8823
$ws.Range("I107").Value = 1032.375
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1032.375
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 887.625
$ws.Range("N107").Value = -5340

$ws.Range("H112").Value = 8286.5
$ws.Range("J112").Value = 4346.1
$ws.Range("L112").Value = 13038.3
$ws.Range("N112").Value = -15254.3

$ws.Range("H122").Value = 1673.8334
$ws.Range("I122").Value = 1780.4546
$ws.Range("K122").Value = 5341.3638
$ws.Range("M122").Value = -2891.3638

$ws.Range("H125").Value = 3984.6365
$ws.Range("I125").Value = 5915.8335
$ws.Range("K125").Value = 53242.5015
$ws.Range("M125").Value = -50782.5015

$ws.Range("H129").Value = 1671.2858
$ws.Range("I129").Value = 924.75
$ws.Range("K129").Value = 2774.25
$ws.Range("M129").Value = 2225.75

$ws.Range("H132").Value = 21452.3
$ws.Range("I132").Value = 12884.412
$ws.Range("K132").Value = 38653.236
$ws.Range("M132").Value = -36123.236

$ws.Range("H137").Value = 5400.2295
$ws.Range("I137").Value = 2023.14
$ws.Range("J137").Value = 20750.637
$ws.Range("K137").Value = 6069.42
$ws.Range("L137").Value = 62251.91099999999
$ws.Range("M137").Value = -3519.42
$ws.Range("N137").Value = -67351.91099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 454.75
$ws.Range("I16").Value = 1500
$ws.Range("J16").Value = 106.333336
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 106.333336
$ws.Range("M16").Value = -1213
$ws.Range("N16").Value = -680.333336

$ws.Range("H32").Value = 5559041
$ws.Range("I32").Value = 8335628.5
$ws.Range("K32").Value = 8335628.5
$ws.Range("M32").Value = -8335341.5

$ws.Range("H45").Value = 3484.1667
$ws.Range("I45").Value = 3484.1667
$ws.Range("K45").Value = 3484.1667
$ws.Range("M45").Value = -3107.1667

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""

$ws.Range("H61").Value = 1025852.25
$ws.Range("I61").Value = 3101.6667
$ws.Range("K61").Value = 3101.6667
$ws.Range("M61").Value = -2889.6667

$ws.Range("H102").Value = 6553
$ws.Range("I102").Value = 6758.3335
$ws.Range("K102").Value = 6758.3335
$ws.Range("M102").Value = -5136.3335

$ws.Range("H110").Value = 6092.609
$ws.Range("I110").Value = 6278.636
$ws.Range("K110").Value = 6278.636
$ws.Range("M110").Value = -4233.636

$ws.Range("H132").Value = 5100122
$ws.Range("I132").Value = 2297.9167
$ws.Range("J132").Value = 9469686
$ws.Range("K132").Value = 6893.750100000001
$ws.Range("L132").Value = 28409058
$ws.Range("M132").Value = -4363.750100000001
$ws.Range("N132").Value = -28414118

$ws.Range("H136").Value = 1025852.25
$ws.Range("I136").Value = 3101.6667
$ws.Range("K136").Value = 9305.000100000001
$ws.Range("M136").Value = -6755.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1323.2916
$ws.Range("I94").Value = 1277.1052
$ws.Range("K94").Value = 1277.1052
$ws.Range("M94").Value = -826.1052

$ws.Range("H99").Value = 14453.423
$ws.Range("I99").Value = 15599.565
$ws.Range("K99").Value = 15599.565
$ws.Range("M99").Value = -14101.565

$ws.Range("H107").Value = 1308.2222
$ws.Range("I107").Value = 1494.5
$ws.Range("J107").Value = 1215.0834
$ws.Range("K107").Value = 1494.5
$ws.Range("L107").Value = 1215.0834
$ws.Range("M107").Value = 425.5
$ws.Range("N107").Value = -5055.0834

$ws.Range("H134").Value = 79301.64999999999
$ws.Range("I134").Value = 155440.86
$ws.Range("K134").Value = 466322.58
$ws.Range("M134").Value = -463787.58

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11076.6875
$ws.Range("I31").Value = 818.44446
$ws.Range("J31").Value = 24265.857
$ws.Range("K31").Value = 818.44446
$ws.Range("L31").Value = 24265.857
$ws.Range("M31").Value = -523.44446
$ws.Range("N31").Value = -24855.857

$ws.Range("H34").Value = 11076.6875
$ws.Range("I34").Value = 818.44446
$ws.Range("J34").Value = 24265.857
$ws.Range("K34").Value = 818.44446
$ws.Range("L34").Value = 24265.857
$ws.Range("M34").Value = -616.44446
$ws.Range("N34").Value = -24669.857

$ws.Range("H58").Value = 17813.5
$ws.Range("I58").Value = 9011.546
$ws.Range("K58").Value = 9011.546
$ws.Range("M58").Value = -8808.546

$ws.Range("H62").Value = 2489.3635
$ws.Range("I62").Value = 2324.8572
$ws.Range("J62").Value = 2777.25
$ws.Range("K62").Value = 2324.8572
$ws.Range("L62").Value = 2777.25
$ws.Range("M62").Value = -1700.8572
$ws.Range("N62").Value = -4025.25

$ws.Range("H65").Value = 2489.3635
$ws.Range("I65").Value = 2324.8572
$ws.Range("J65").Value = 2777.25
$ws.Range("K65").Value = 11624.286
$ws.Range("L65").Value = 13886.25
$ws.Range("M65").Value = -8504.286
$ws.Range("N65").Value = -20126.25

$ws.Range("H120").Value = 77775
$ws.Range("J120").Value = 77775
$ws.Range("L120").Value = 77775
$ws.Range("N120").Value = -85033

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""

$ws.Range("H132").Value = 31197016
$ws.Range("I132").Value = 2593.48
$ws.Range("K132").Value = 7780.440000000001
$ws.Range("M132").Value = -5250.440000000001

$ws.Range("H136").Value = 17813.5
$ws.Range("I136").Value = 9011.546
$ws.Range("K136").Value = 27034.638
$ws.Range("M136").Value = -24484.638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 6097.278
$ws.Range("I7").Value = 125.166664
$ws.Range("K7").Value = 375.499992
$ws.Range("M7").Value = -263.499992

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").Value = ""

$ws.Range("H92").Value = 320.2
$ws.Range("I92").Value = 233.66667
$ws.Range("K92").Value = 701.00001
$ws.Range("M92").Value = 546.99999

$ws.Range("H104").Value = 4118559.8
$ws.Range("I104").Value = 5000
$ws.Range("J104").Value = 12345679
$ws.Range("K104").Value = 15000
$ws.Range("L104").Value = 37037037
$ws.Range("M104").Value = -12379
$ws.Range("N104").Value = -37042279

$ws.Range("H112").Value = 4508.6665
$ws.Range("J112").Value = 7500
$ws.Range("L112").Value = 22500
$ws.Range("N112").Value = -24716

$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").Value = ""

$ws.Range("H122").Value = 12662427
$ws.Range("J122").Value = 2365959.5
$ws.Range("L122").Value = 21293635.5
$ws.Range("N122").Value = -21298535.5

$ws.Range("H125").Value = 30
$ws.Range("I125").Value = 30
$ws.Range("K125").Value = 90
$ws.Range("M125").Value = 4830

$ws.Range("H129").Value = 16668612
$ws.Range("I129").Value = 20002190
$ws.Range("J129").Value = 725
$ws.Range("K129").Value = 60006570
$ws.Range("L129").Value = 2175
$ws.Range("M129").Value = -60001570
$ws.Range("N129").Value = -12175

$ws.Range("H131").Value = 3478.9082
$ws.Range("I131").Value = 12382.5
$ws.Range("J131").Value = 2687.4778
$ws.Range("K131").Value = 37147.5
$ws.Range("L131").Value = 8062.4334
$ws.Range("M131").Value = -32107.5
$ws.Range("N131").Value = -18142.4334

$ws.Range("H137").Value = 10707.5
$ws.Range("I137").Value = 6999
$ws.Range("J137").Value = 11943.667
$ws.Range("K137").Value = 20997
$ws.Range("L137").Value = 35831.001
$ws.Range("M137").Value = -15897
$ws.Range("N137").Value = -46031.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 2916.6667
$ws.Range("J6").Value = 2916.6667
$ws.Range("L6").Value = 2916.6667
$ws.Range("N6").Value = -3142.6667

$ws.Range("H7").Value = 556088.9
$ws.Range("J7").Value = 1667266.6
$ws.Range("L7").Value = 1667266.6
$ws.Range("N7").Value = -1667490.6

$ws.Range("H8").Value = 556088.9
$ws.Range("J8").Value = 1667266.6
$ws.Range("L8").Value = 1667266.6
$ws.Range("N8").Value = -1667544.6

$ws.Range("H16").Value = 2916.6667
$ws.Range("J16").Value = 2916.6667
$ws.Range("L16").Value = 2916.6667
$ws.Range("N16").Value = -3416.6667

$ws.Range("H97").Value = 1370.2162
$ws.Range("I97").Value = 1066.0625
$ws.Range("K97").Value = 1066.0625
$ws.Range("M97").Value = -570.0625

$ws.Range("H102").Value = 6241.25
$ws.Range("I102").Value = 7571.6665
$ws.Range("K102").Value = 7571.6665
$ws.Range("M102").Value = -5949.6665

$ws.Range("H122").Value = 2855.9333
$ws.Range("I122").Value = 3181.25
$ws.Range("K122").Value = 9543.75
$ws.Range("M122").Value = -7093.75

$ws.Range("H132").Value = 634573.7
$ws.Range("I132").Value = 3889.9412
$ws.Range("J132").Value = 1706736.1
$ws.Range("K132").Value = 11669.8236
$ws.Range("L132").Value = 5120208.300000001
$ws.Range("M132").Value = -9139.8236
$ws.Range("N132").Value = -5125268.300000001

$ws.Range("H141").Value = 112594
$ws.Range("J141").Value = 112594
$ws.Range("L141").Value = 112594
$ws.Range("N141").Value = -122954

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6151.7393
$ws.Range("I7").Value = 6535
$ws.Range("J7").Value = 5555.5557
$ws.Range("K7").Value = 6535
$ws.Range("L7").Value = 5555.5557
$ws.Range("M7").Value = -6423
$ws.Range("N7").Value = -5779.5557

$ws.Range("H40").Value = 3512.9048
$ws.Range("I40").Value = 2485.6875
$ws.Range("K40").Value = 2485.6875
$ws.Range("M40").Value = -2349.6875

$ws.Range("H68").Value = 10408.482
$ws.Range("I68").Value = 8593.214
$ws.Range("K68").Value = 8593.214
$ws.Range("M68").Value = -7844.214

$ws.Range("H71").Value = 10408.482
$ws.Range("I71").Value = 8593.214
$ws.Range("K71").Value = 42966.07
$ws.Range("M71").Value = -39222.07

$ws.Range("H82").Value = 2407
$ws.Range("I82").Value = 2623.3076
$ws.Range("J82").Value = 2055.5
$ws.Range("K82").Value = 2623.3076
$ws.Range("L82").Value = 2055.5
$ws.Range("M82").Value = -2262.3076
$ws.Range("N82").Value = -2777.5

$ws.Range("H85").Value = 2407
$ws.Range("I85").Value = 2623.3076
$ws.Range("J85").Value = 2055.5
$ws.Range("K85").Value = 2623.3076
$ws.Range("L85").Value = 2055.5
$ws.Range("M85").Value = -1375.3076
$ws.Range("N85").Value = -4551.5

$ws.Range("H96").Value = 16666.666
$ws.Range("J96").Value = 16666.666
$ws.Range("L96").Value = 16666.666
$ws.Range("N96").Value = -22158.666

$ws.Range("H122").Value = 5950.1377
$ws.Range("I122").Value = 4680.0835
$ws.Range("J122").Value = 6846.647
$ws.Range("K122").Value = 14040.2505
$ws.Range("L122").Value = 20539.941
$ws.Range("M122").Value = -11590.2505
$ws.Range("N122").Value = -25439.941

$ws.Range("H126").Value = 6151.7393
$ws.Range("I126").Value = 6535
$ws.Range("J126").Value = 5555.5557
$ws.Range("K126").Value = 19605
$ws.Range("L126").Value = 16666.6671
$ws.Range("M126").Value = -17135
$ws.Range("N126").Value = -21606.6671

$ws.Range("H132").Value = 2337343.8
$ws.Range("I132").Value = 3190
$ws.Range("K132").Value = 9570
$ws.Range("M132").Value = -7040

$ws.Range("H136").Value = 973433.4
$ws.Range("I136").Value = 16988.715
$ws.Range("J136").Value = 1678182.1
$ws.Range("K136").Value = 50966.145
$ws.Range("L136").Value = 5034546.300000001
$ws.Range("M136").Value = -48416.145
$ws.Range("N136").Value = -5039646.300000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 219184.5
$ws.Range("I9").Value = 219184.5
$ws.Range("K9").Value = 219184.5
$ws.Range("M9").Value = -219044.5

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").Value = ""

$ws.Range("H14").Value = 4868.3335
$ws.Range("J14").Value = 4868.3335
$ws.Range("L14").Value = 4868.3335
$ws.Range("N14").Value = -5204.3335

$ws.Range("H33").Value = 46092
$ws.Range("J33").Value = 46092
$ws.Range("L33").Value = 46092
$ws.Range("N33").Value = -46592

$ws.Range("H36").Value = 46092
$ws.Range("J36").Value = 46092
$ws.Range("L36").Value = 46092
$ws.Range("N36").Value = -46592

$ws.Range("H54").Value = 44943.668
$ws.Range("I54").Value = 20570
$ws.Range("J54").Value = 57130.5
$ws.Range("K54").Value = 20570
$ws.Range("L54").Value = 57130.5
$ws.Range("M54").Value = -20050
$ws.Range("N54").Value = -58170.5

$ws.Range("H122").Value = 3438.2896
$ws.Range("I122").Value = 2010.1923
$ws.Range("K122").Value = 6030.5769
$ws.Range("M122").Value = -3580.5769

$ws.Range("H132").Value = 519886.66
$ws.Range("I132").Value = 2758.2856
$ws.Range("K132").Value = 8274.856800000001
$ws.Range("M132").Value = -5744.856800000001

$ws.Range("H136").Value = 579922.4399999999
$ws.Range("I136").Value = 2463.3845
$ws.Range("K136").Value = 7390.1535
$ws.Range("M136").Value = -4840.1535
